# Update the three-digit-by-one-digit multiplication answers in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("387×3=1161", "721×2=1442"),
    @("452×8=3616", "680×7=4760"),
    @("661×6=3966", "568×9=5112"),
    @("447×3=1341", "554×8=4432"),
    @("134×6=804",  "185×3=555"),
    @("456×4=1824", "995×7=6965"),
    @("344×2=688",  "850×8=6800"),
    @("393×4=1572", "941×7=6587"),
    @("696×9=6264", "567×2=1134"),
    @("811×4=3244", "709×9=6381"),
    @("733×2=1466", "532×8=4256"),
    @("935×6=5610", "724×7=5068"),
    @("593×8=4744", "664×6=3984"),
    @("875×9=7875", "808×8=6464"),
    @("318×7=2226", "110×7=770"),
    @("532×2=1064", "853×4=3412"),
    @("436×4=1744", "492×5=2460"),
    @("675×2=1350", "461×9=4149"),
    @("625×2=1250", "222×8=1776"),
    @("513×4=2052", "118×4=472"),
    @("284×7=1988", "687×6=4122"),
    @("925×3=2775", "376×3=1128"),
    @("525×3=1575", "448×2=896"),
    @("307×3=921",  "150×4=600"),
    @("277×9=2493", "486×8=3888")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
